{"js": "// Apply the same content edits as the target diff:\n//  1. Collapse the \">>> your stuff after this line >>>\" paragraph\n//     (previously split across several runs around proofErr markers)\n//     into a single plain run.\n//  2. Change \"Ben changing things up!\" to \"Suman Rijal File changed.\"\n//     and leave a \"_GoBack\" bookmark right after the new text (Word\n//     drops this bookmark at the point of the most recent edit).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraphs by their current text rather than hard-coded\n// indexes, so the script is resilient to minor document variations.\nlet markerPara = null;\nlet benPara = null;\nfor (const p of paragraphs.items) {\n  const t = p.text;\n  if (t.indexOf(\">>\") !== -1 && t.indexOf(\"your\") !== -1 && t.indexOf(\"stuff after this line\") !== -1) {\n    markerPara = p;\n  } else if (t.indexOf(\"Ben changing things up!\") !== -1) {\n    benPara = p;\n  }\n}\n\nif (markerPara) {\n  // Re-write as a single, un-split run: \">>>  your stuff after this line >>>\"\n  markerPara.insertText(\">>>  your stuff after this line >>>\", Word.InsertLocation.replace);\n}\n\nif (benPara) {\n  benPara.insertText(\"Suman Rijal File changed.\", Word.InsertLocation.replace);\n  await context.sync();\n\n  // Drop a \"_GoBack\" bookmark right after the replaced text, mirroring\n  // the bookmark Word leaves behind at the site of the last edit.\n  const endRange = benPara.getRange(Word.RangeLocation.end);\n  endRange.insertBookmark(\"_GoBack\");\n}\n\nawait context.sync();\n", "ps1": "# Apply the same content edits as the target diff:\n#  1. Collapse the \">>> your stuff after this line >>>\" paragraph\n#     (previously split across several runs around proofErr markers)\n#     into a single plain run.\n#  2. Change \"Ben changing things up!\" to \"Suman Rijal File changed.\"\n#     and leave a \"_GoBack\" bookmark right after the new text (Word\n#     drops this bookmark at the point of the most recent edit).\n\n$d = $word.ActiveDocument\n\n# --- 1. Re-write the \">>> ... >>>\" marker paragraph as a single run ---\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $candidate = $d.Paragraphs.Item($i)\n    $candidateText = $candidate.Range.Text\n    if ($candidateText -like \"*your*\" -and $candidateText -like \"*stuff after this line*\") {\n        $markerRange = $candidate.Range\n        $markerRange.MoveEnd(1, -1) | Out-Null\n        # Clear then re-insert so Word rebuilds a single fresh run instead\n        # of leaving the text spread across the old proofErr-split runs.\n        $markerRange.Text = \"\"\n        $markerRange.InsertAfter(\">>>  your stuff after this line >>>\")\n        break\n    }\n}\n\n# --- 2. Replace \"Ben changing things up!\" and drop a _GoBack bookmark ---\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $candidate = $d.Paragraphs.Item($i)\n    if ($candidate.Range.Text -like \"*Ben changing things up!*\") {\n        $benRange = $candidate.Range\n        $benRange.MoveEnd(1, -1) | Out-Null\n        $benRange.Text = \"Suman Rijal File changed.\"\n\n        # Leave a one-character sentinel after the text so we can anchor a\n        # zero-width bookmark exactly at end-of-text, then remove it again.\n        $para = $d.Paragraphs.Item($i)\n        $endRange = $para.Range\n        $endRange.MoveEnd(1, -1) | Out-Null\n        $endRange.InsertAfter(\"X\")\n\n        $para2 = $d.Paragraphs.Item($i)\n        $fullRange = $para2.Range\n        $fullRange.MoveEnd(1, -1) | Out-Null\n        $bmRange = $fullRange.Duplicate\n        $bmRange.Start = $fullRange.End - 1\n        $bmRange.End = $fullRange.End - 1\n        $d.Bookmarks.Add(\"_GoBack\", $bmRange) | Out-Null\n\n        $para3 = $d.Paragraphs.Item($i)\n        $trimRange = $para3.Range\n        $trimRange.MoveEnd(1, -1) | Out-Null\n        $sentinelRange = $trimRange.Duplicate\n        $sentinelRange.Start = $trimRange.End - 1\n        $sentinelRange.End = $trimRange.End\n        $sentinelRange.Text = \"\"\n        break\n    }\n}\n"}
